$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Max Incoming Interactions" / "Max Outgoing Interactions" headers
# (order of cluster in heatmaps)
$ws.Range("C1").Value = "Max Outgoing Interactions"
$ws.Range("E1").Value = "Max Incoming Interactions"
